$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the transfer-result messages in column F (rows 3 and 5)
$ws.Range("F3").Value = "successfully transferred"
$ws.Range("F5").Value = "successfully transferred"

# Populate the new "Prueba OK" markers in column G for rows 3-5
$ws.Range("G3").Value = "Prueba OK"
$ws.Range("G4").Value = "Prueba OK"
$ws.Range("G5").Value = "Prueba OK"

# Update the active cell selection to F5
$ws.Range("F5").Select()
